# Removing less than USD 5 price from extrapolation calibration because it is just noise.
# Updates the recalculated ABSM1_RN / M1_RN / CM2_RN / CMN3_RN / CMN4_RN columns (D:H)
# for the rows whose calibration results changed once the sub-$5 price point was excluded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = @{ D = 112496.5643774469;  E = -0.02751461436739952; F = 0.1908283209229917;  G = -1.630757786286523;  H = 13.94330982624902 }
    5  = @{ D = 113257.7075729683;  E = -0.0260929860268405;  F = 0.1976417838309602;  G = -0.9526286791957003; H = 8.011327814748469 }
    7  = @{ D = 115225.6038564603;  E = -0.05720997373000878; F = 0.357941707382721;   G = -2.143755865867489;  H = 10.360678328631 }
    10 = @{ D = 118504.1672356246;  E = -0.1198421171135392;  F = 0.4463146263340288;  G = -1.87982023229993;   H = 9.368141195375056 }
    11 = @{ D = 120171.7317044994;  E = -0.197366911094061;   F = 0.7694093643850638;  G = -2.508656671379411;  H = 11.85477798019164 }
    16 = @{ D = 110850.5560671159;  E = -0.1240136578688794;  F = 0.2112073882646306;  G = -1.23600031829959;   H = 6.953670483699389 }
    19 = @{ D = 111852.8801667728;  E = -0.02543629067531163; F = 0.1502538349772886;  G = -0.4126904493292714; H = 5.976168422060785 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
